$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7 and 8, column H: "false" -> "true" (reuse exact existing shared string "true" from H6)
$ws.Range("H6").Copy()
$ws.Range("H7").PasteSpecial(-4163)
$ws.Range("H6").Copy()
$ws.Range("H8").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Add new column I ("id") with header formatted like the other headers (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "id"

# Fill I2:I14 with sequential numeric ids 0..12
$ids = 0,1,2,3,4,5,6,7,8,9,10,11,12
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $ids[$i]
}

# Update the active selection to match the edited workbook
$ws.Range("H8").Select()
